# Applies the "added subrogation expense and cause of loss" edit to
# ClaimMapping.xlsx:
#   - Mapping sheet: fix "Casue Of Loss" typo, expand "eBao" -> "eBao Class",
#     and append three new VesselRiskinfo mapping rows (20-22).
#   - DefaultValues sheet: update the default values used for
#     DefaultUsernameForSuClaimOwner and CauseofLossDropdown.
#   - Leave the workbook with the DefaultValues tab active/selected (B3),
#     matching the last place the author clicked before saving, with the
#     Mapping sheet's own selection parked at B22.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Mapping")
$ws2 = $wb.Worksheets.Item("DefaultValues")

# --- Mapping sheet fixes -----------------------------------------------
$ws1.Range("C17").Value = "Cause Of Loss"
$ws1.Range("C19").Value = "eBao Class"

# --- Mapping sheet: new VesselRiskinfo rows -----------------------------
$ws1.Range("B20").Value = "VesselPolicyNumber"
$ws1.Range("C20").Value = "policy"
$ws1.Range("D20").Value = "From VesselRiskinfo excel"

$ws1.Range("B21").Value = "VesselClaimNumber"
$ws1.Range("C21").Value = "claim"
$ws1.Range("D21").Value = "From VesselRiskinfo excel"

$ws1.Range("B22").Value = "VesselName"
$ws1.Range("C22").Value = "Vessel Name"
$ws1.Range("D22").Value = "From VesselRiskinfo excel"

# --- DefaultValues sheet updates ----------------------------------------
$ws2.Range("B2").Value = "Project One"
$ws2.Range("B3").Value = "Others"

# --- Selections / active tab: DefaultValues ends up the active sheet,
# selected at B3; Mapping's own lingering selection is B22.
[void]$ws1.Range("B22").Select()
$ws2.Activate() | Out-Null
[void]$ws2.Range("B3").Select()
